# Commit: "edit the sheet name test"
# The "RSPP see note" tab is renamed to "RSPP" (RSPP was acquired by CXO,
# so the "see note" caveat on that tab name is no longer needed).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("RSPP see note")
$ws.Name = "RSPP"

# While testing the rename, the "CXO see note" tab ends up the active/
# selected sheet at save time.
$wb.Worksheets.Item("CXO see note").Activate()
